$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.310.33"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.431.37"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.04"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.76"
$ws.Range("E6").Value = "  +1.22%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.533"
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.430.62"
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("E12").Value = "  -1.26%  "
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.65"
$ws.Range("E14").Value = "  +3.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000176"
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.102.30"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.424.17"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.29"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "324.85"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("E22").Value = "  -1.22%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.42"
$ws.Range("E24").Value = "  +2.45%  "
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("E26").Value = "  -3.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "553.81"
$ws.Range("E27").Value = "  -4.23%  "
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.29"
$ws.Range("E31").Value = "  +0.60%  "
$ws.Range("E32").Value = "  -1.19%  "
$ws.Range("E33").Value = "  -1.29%  "
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("E35").Value = "  -1.15%  "
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.82"
$ws.Range("E37").Value = "  +1.63%  "
$ws.Range("E38").Value = "  -0.93%  "
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.72"
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "150.10"
$ws.Range("E41").Value = "  -1.30%  "
$ws.Range("E42").Value = "  -0.76%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "147.93"
$ws.Range("E45").Value = "  -1.07%  "
$ws.Range("E46").Value = "  +0.50%  "
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.29"
$ws.Range("E48").Value = "  +0.87%  "
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0925"
$ws.Range("E50").Value = "  +0.77%  "
$ws.Range("E51").Value = "  +0.89%  "
